$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new rows 31-45 with CrM28-CrM38 and CaM0-CaM3 test cases
# Row 31
$ws.Range("D31").Value = "CrM28"
$ws.Range("E31").Value = "Crear Modelo"
$ws.Range("F31").Value = "CrM23||CrM24"
$ws.Range("G31").Value = "1-Se prueba a redactar una descripción del modelo"
$ws.Range("H31").Value = "La descripción se muestra correctamente"
$ws.Range("I31").Value = "✅"

# Row 32
$ws.Range("D32").Value = "CrM29"
$ws.Range("E32").Value = "Guardar Modelo"
$ws.Range("F32").Value = "CrM23||CrM24"
$ws.Range("G32").Value = "1-Se presiona el botón `"Guardar Modelo`" "
$ws.Range("H32").Value = "Ventana Emergete para seleccionar carpeta"
$ws.Range("I32").Value = "✅"

# Row 33
$ws.Range("D33").Value = "CrM30"
$ws.Range("E33").Value = "Guardar Modelo"
$ws.Range("F33").Value = "CrM29"
$ws.Range("G33").Value = "1-Se selecciona una carpeta sin permisos (Ejemplo: System32)"
$ws.Range("H33").Value = "Notificación de Error `"Permission denied`""
$ws.Range("I33").Value = "✅"

# Row 34
$ws.Range("D34").Value = "CrM31"
$ws.Range("E34").Value = "Guardar Modelo"
$ws.Range("F34").Value = "CrM29"
$ws.Range("G34").Value = "1-Se selecciona un disco sin espacio suficiente"
$ws.Range("H34").Value = "Notifcación de Error"
$ws.Range("I34").Value = "❓"

# Row 35
$ws.Range("D35").Value = "CrM32"
$ws.Range("E35").Value = "Guardar Modelo"
$ws.Range("F35").Value = "CrM29"
$ws.Range("G35").Value = "1-Se selecciona una carpeta válida"
$ws.Range("H35").Value = "Notificación de Éxito"
$ws.Range("I35").Value = "✅"

# Row 36
$ws.Range("D36").Value = "CrM33"
$ws.Range("E36").Value = "Guardar Modelo"
$ws.Range("F36").Value = "CrM29"
$ws.Range("G36").Value = "1-Se acepta o cierra la notficación"
$ws.Range("H36").Value = "Se Guarda el modelo/descripción"
$ws.Range("I36").Value = "✅"

# Row 37
$ws.Range("D37").Value = "CrM34"
$ws.Range("E37").Value = "Predicción del Modelo"
$ws.Range("F37").Value = "CrM33||CrM23||CrM24"
$ws.Range("G37").Value = "1-Se comprueba el panel `"Predicción con modelo`""
$ws.Range("H37").Value = "Panel correcto para la predicción"
$ws.Range("I37").Value = "✅"

# Row 38
$ws.Range("D38").Value = "CrM35"
$ws.Range("E38").Value = "Predicción del Modelo"
$ws.Range("F38").Value = "CrM34"
$ws.Range("G38").Value = "1-Se comprueba que la entrada/salida son las respectivas columnas seleccionadas"
$ws.Range("H38").Value = "Panel dividido correctamente"
$ws.Range("I38").Value = "✅"

# Row 39
$ws.Range("D39").Value = "CrM36"
$ws.Range("E39").Value = "Predicción del Modelo"
$ws.Range("F39").Value = "CrM34"
$ws.Range("G39").Value = "1-No se escribe una entrada o se escribe un valor no numérico y se predice"
$ws.Range("H39").Value = "Notificación de Error `"Valores no numéricos`""
$ws.Range("I39").Value = "✅"

# Row 40
$ws.Range("D40").Value = "CrM37"
$ws.Range("E40").Value = "Predicción del Modelo"
$ws.Range("F40").Value = "CrM34"
$ws.Range("G40").Value = "1-Se escribe un valor numérico en la entrada y se predice"
$ws.Range("H40").Value = "Se predice correctamente la salida"
$ws.Range("I40").Value = "✅"

# Row 41
$ws.Range("D41").Value = "CrM38"
$ws.Range("E41").Value = "Cargar Datos"
$ws.Range("F41").Value = "CrM(1||2||…..||37)"
$ws.Range("G41").Value = "1-Se carga otro archivo"
$ws.Range("H41").Value = "Se reinician todos los paneles"
$ws.Range("I41").Value = "✅"

# Row 42
$ws.Range("D42").Value = "CaM0"
$ws.Range("E42").Value = "Cargar Modelo"
$ws.Range("F42").Value = "IN0"
$ws.Range("G42").Value = "1-Se presiona el botón `"Cargar Modelo`" arriba a la izquierda"
$ws.Range("H42").Value = "Nueva pestaña con paneles para cargar modelo"
$ws.Range("I42").Value = "✅"

# Row 43
$ws.Range("D43").Value = "CaM1"
$ws.Range("E43").Value = "Cargar Modelo"
$ws.Range("F43").Value = "CaM0"
$ws.Range("G43").Value = "1-Se presiona el botón`"Cargar Modelo`" del panel `"Cargar Modelo Guardado`""
$ws.Range("H43").Value = "Nueva ventana para selección de archivo"
$ws.Range("I43").Value = "✅"

# Row 44
$ws.Range("D44").Value = "CaM2"
$ws.Range("E44").Value = "Cargar Modelo"
$ws.Range("F44").Value = "CaM1"
$ws.Range("G44").Value = "1-Se exploran archivos"
$ws.Range("H44").Value = "Solo se aceptan arhcivos Joblib"
$ws.Range("I44").Value = "✅"

# Row 45
$ws.Range("D45").Value = "CaM3"
$ws.Range("E45").Value = "Cargar Modelo"
$ws.Range("F45").Value = "CaM2"
$ws.Range("G45").Value = "1-Se selecciona un archivo corrupto"
$ws.Range("H45").Value = "Notificación de Error `"Error al cargar el modelo`""
$ws.Range("I45").Value = "✅"

# Move the bottom-border style (previously on G38) to G37, matching the new table layout
$ws.Range("G37").Borders.Item(9).Color = 3381759
$ws.Range("G37").Borders.Item(9).LineStyle = 1
$ws.Range("G38").Borders.Item(9).LineStyle = -4142

# Resize Tabla1 to cover the newly added rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D1:I46"))

# Update the active selection to reflect the area that was being edited
$ws.Range("K38").Select()
